# Refresh the "cryptos" price/volume snapshot for this run.
# Each entry is a target cell plus its new literal text. Cells whose new
# text reads as a plain number (e.g. "518.47") are written with a leading
# apostrophe so Excel stores them as text (matching column D/E's existing
# inline-string convention) instead of silently coercing them to numbers;
# ClearFormats() afterwards drops the quote-prefix style Excel applies for
# that trick so the cell keeps its original (style-less) formatting.
$updates = @(
    @{ Cell = 'D2'; Text = '57.976.45' }
    @{ Cell = 'E2'; Text = '  +2.78%  ' }
    @{ Cell = 'D3'; Text = '3.051.41' }
    @{ Cell = 'E3'; Text = '  +2.37%  ' }
    @{ Cell = 'E4'; Text = '  +0.03%  ' }
    @{ Cell = 'D5'; Text = '518.47' }
    @{ Cell = 'E5'; Text = '  +2.95%  ' }
    @{ Cell = 'D6'; Text = '141.90' }
    @{ Cell = 'E6'; Text = '  +5.37%  ' }
    @{ Cell = 'D7'; Text = '1.00' }
    @{ Cell = 'E7'; Text = '  -0.03%  ' }
    @{ Cell = 'E8'; Text = '  +4.07%  ' }
    @{ Cell = 'E9'; Text = '  +2.95%  ' }
    @{ Cell = 'E10'; Text = '  +5.47%  ' }
    @{ Cell = 'D11'; Text = '0.370' }
    @{ Cell = 'E11'; Text = '  +5.55%  ' }
    @{ Cell = 'D12'; Text = '3.575.32' }
    @{ Cell = 'E12'; Text = '  +2.44%  ' }
    @{ Cell = 'E13'; Text = '  +2.23%  ' }
    @{ Cell = 'D14'; Text = '26.74' }
    @{ Cell = 'E14'; Text = '  +6.67%  ' }
    @{ Cell = 'E15'; Text = '  +12.13%  ' }
    @{ Cell = 'D16'; Text = '58.006.06' }
    @{ Cell = 'E16'; Text = '  +2.84%  ' }
    @{ Cell = 'D17'; Text = '6.25' }
    @{ Cell = 'E17'; Text = '  +10.34%  ' }
    @{ Cell = 'D18'; Text = '3.052.79' }
    @{ Cell = 'D19'; Text = '13.04' }
    @{ Cell = 'E19'; Text = '  +5.92%  ' }
    @{ Cell = 'D20'; Text = '8.10' }
    @{ Cell = 'E20'; Text = '  +4.64%  ' }
    @{ Cell = 'D21'; Text = '338.52' }
    @{ Cell = 'E21'; Text = '  +4.60%  ' }
    @{ Cell = 'E22'; Text = '  +0.04%  ' }
    @{ Cell = 'D23'; Text = '5.77' }
    @{ Cell = 'E23'; Text = '  +1.35%  ' }
    @{ Cell = 'D24'; Text = '0.502' }
    @{ Cell = 'E24'; Text = '  +6.89%  ' }
    @{ Cell = 'D25'; Text = '65.11' }
    @{ Cell = 'E25'; Text = '  +5.24%  ' }
    @{ Cell = 'E26'; Text = '  +3.94%  ' }
    @{ Cell = 'D27'; Text = '0.0₃0956' }
    @{ Cell = 'E27'; Text = '  +7.55%  ' }
    @{ Cell = 'D28'; Text = '1.00' }
    @{ Cell = 'E28'; Text = '  +0.73%  ' }
    @{ Cell = 'D29'; Text = '6.92' }
    @{ Cell = 'E29'; Text = '  +6.16%  ' }
    @{ Cell = 'D30'; Text = '7.59' }
    @{ Cell = 'E30'; Text = '  +12.22%  ' }
    @{ Cell = 'B31'; Text = 'Fetch.AI' }
    @{ Cell = 'C31'; Text = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet' }
    @{ Cell = 'D31'; Text = '1.24' }
    @{ Cell = 'E31'; Text = '  +4.04%  ' }
    @{ Cell = 'B32'; Text = 'PancakeSwap' }
    @{ Cell = 'C32'; Text = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D32'; Text = '1.83' }
    @{ Cell = 'E32'; Text = '  +4.60%  ' }
    @{ Cell = 'E33'; Text = '  +3.04%  ' }
    @{ Cell = 'D34'; Text = '4.78' }
    @{ Cell = 'E34'; Text = '  +7.31%  ' }
    @{ Cell = 'D35'; Text = '156.42' }
    @{ Cell = 'E35'; Text = '  -0.75%  ' }
    @{ Cell = 'D36'; Text = '5.93' }
    @{ Cell = 'E36'; Text = '  +7.05%  ' }
    @{ Cell = 'D37'; Text = '1.30' }
    @{ Cell = 'E37'; Text = '  +2.19%  ' }
    @{ Cell = 'D38'; Text = '25.31' }
    @{ Cell = 'E38'; Text = '  +10.33%  ' }
    @{ Cell = 'D39'; Text = '0.0693' }
    @{ Cell = 'E39'; Text = '  +2.89%  ' }
    @{ Cell = 'D40'; Text = '3.086.65' }
    @{ Cell = 'E40'; Text = '  +2.37%  ' }
    @{ Cell = 'D41'; Text = '37.78' }
    @{ Cell = 'E41'; Text = '  +4.19%  ' }
    @{ Cell = 'D42'; Text = '3.91' }
    @{ Cell = 'E42'; Text = '  +10.13%  ' }
    @{ Cell = 'E43'; Text = '  +0.07%  ' }
    @{ Cell = 'E44'; Text = '  +4.03%  ' }
    @{ Cell = 'D45'; Text = '2.331.22' }
    @{ Cell = 'E45'; Text = '  +3.71%  ' }
    @{ Cell = 'E46'; Text = '  +4.34%  ' }
    @{ Cell = 'E47'; Text = '  +2.58%  ' }
    @{ Cell = 'D48'; Text = '6.08' }
    @{ Cell = 'E48'; Text = '  +5.51%  ' }
    @{ Cell = 'D49'; Text = '0.0242' }
    @{ Cell = 'E49'; Text = '  +2.93%  ' }
    @{ Cell = 'D50'; Text = '19.84' }
    @{ Cell = 'E50'; Text = '  +4.80%  ' }
    @{ Cell = 'E51'; Text = '  -3.17%  ' }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    $text = $update.Text
    $looksNumeric = $text -match '^\s*[+-]?(\d+\.?\d*|\.\d+)\s*$'
    if ($looksNumeric) {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
    $range.ClearFormats()
}
